$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("G2").Value = "9840032111"
    $ws.Range("AF2").Value = "9840026220"
    $ws.Range("AV2").Value = "9840087286"
    $ws.Range("AZ2").Value = "9840054625"
}

$sheet1 = $wb.Worksheets.Item(1)
$sheet1.Range("O2").Value = "06-01-2025"
$sheet1.Range("Q2").Value = "09-01-2025 05:00:00 PM"
$sheet1.Range("AD2").Value = "06-01-2025"
$sheet1.Range("BB2").Value = "CT: Mon, Jan 06, 2025 at 8:52 PM"

$sheet2 = $wb.Worksheets.Item(2)
$sheet2.Range("O2").Value = "06-01-2025"
$sheet2.Range("Q2").Value = "09-01-2025 05:00:00 PM"
$sheet2.Range("AD2").Value = "06-01-2025"
$sheet2.Range("BB2").Value = "CT: Mon, Jan 06, 2025 at 9:01 PM"
